# Update CDA Logical model for ST.r2b
# - Rename the "Include from ActMood" sheet to "Include #0"
# - Bump Version / Date on the Metadata sheet
# - Insert a new "Jurisdiction" property row (with an empty value) right
#   after "Contact", pushing Description/Purpose/Copyright/Immutable down
#   by one row.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)
$include = $wb.Worksheets.Item(2)

# --- rename the include sheet -------------------------------------------
$include.Name = "Include #0"

# --- bump the Version / Date metadata values ----------------------------
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- insert the new Jurisdiction row at row 11 --------------------------
# This shifts the existing rows 11-14 (Description, Purpose, Copyright,
# Immutable) down to rows 12-15, carrying their formatting with them.
$meta.Rows.Item(11).Insert()

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Copy the row-10 (Contact) formatting onto the freshly inserted row 11 so
# it matches the rest of the table (border/fill/alignment) instead of the
# blank default style Insert() leaves behind.
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
